$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 2 4 "321.29"
Set-TextValue 2 5 "-3.08%"
Set-TextValue 3 4 "42.80"
Set-TextValue 3 5 "-5.91%"
Set-TextValue 4 4 "5.193"
Set-TextValue 4 5 "-5.40%"
Set-TextValue 5 5 "-3.28%"
Set-TextValue 6 4 "4.326"
Set-TextValue 6 5 "-2.69%"
Set-TextValue 7 4 "1.774"
Set-TextValue 7 5 "-14.03%"
Set-TextValue 8 4 "0.9500"
Set-TextValue 8 5 "-3.39%"
Set-TextValue 9 4 "0.1121"
Set-TextValue 9 5 "-3.37%"
Set-TextValue 10 4 "0.1886"
Set-TextValue 10 5 "-1.50%"
Set-TextValue 11 4 "0.09367"
Set-TextValue 11 5 "-3.95%"
Set-TextValue 12 4 "0.04650"
Set-TextValue 12 5 "-1.33%"
Set-TextValue 13 4 "7.448"
Set-TextValue 13 5 "-21.16%"
Set-TextValue 14 5 "-0.37%"
Set-TextValue 15 4 "0.001293"
Set-TextValue 15 5 "-0.62%"
Set-TextValue 16 4 "0.005706"
Set-TextValue 16 5 "-4.13%"
Set-TextValue 17 4 "3.357"
Set-TextValue 17 5 "-0.93%"
Set-TextValue 18 4 "2.589"
Set-TextValue 18 5 "2.18%"
Set-TextValue 19 5 "0.32%"
Set-TextValue 20 4 "0.1389"
Set-TextValue 20 5 "2.53%"
Set-TextValue 22 4 "0.04175"
Set-TextValue 22 5 "0.84%"
Set-TextValue 23 4 "0.001252"
Set-TextValue 23 5 "-3.88%"
Set-TextValue 24 4 "0.004303"
Set-TextValue 24 5 "-3.65%"
Set-TextValue 25 4 "0.0001222"
Set-TextValue 25 5 "-6.34%"
Set-TextValue 38 4 "0.02668"
Set-TextValue 38 5 "-4.26%"
Set-TextValue 39 4 "0.05638"
Set-TextValue 39 5 "-1.17%"
Set-TextValue 40 4 "0.008116"
Set-TextValue 40 5 "3.74%"
Set-TextValue 41 4 "0.1405"
Set-TextValue 41 5 "-1.92%"
Set-TextValue 42 4 "0.006523"
Set-TextValue 42 5 "-10.25%"
Set-TextValue 43 4 "0.002120"
Set-TextValue 43 5 "-0.10%"
Set-TextValue 44 4 "0.007688"
Set-TextValue 44 5 "-3.20%"
Set-TextValue 45 4 "0.3488"
Set-TextValue 45 5 "2.59%"
Set-TextValue 46 4 "0.00006775"
Set-TextValue 46 5 "-3.12%"
Set-TextValue 47 4 "0.00000000750"
Set-TextValue 47 5 "-0.31%"
Set-TextValue 48 4 "0.003369"
Set-TextValue 48 5 "-2.69%"
Set-TextValue 49 5 "15.88%"
Set-TextValue 50 5 "-0.31%"
Set-TextValue 51 5 "-0.31%"
